# Fixed some bugs in stats
# The data rows (A2:F25) got reshuffled into a different (corrected) row order.
# Capture the original values first, then write them back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (both referring to the ORIGINAL layout)
$mapping = @{
    2  = 3
    3  = 5
    4  = 4
    5  = 6
    6  = 2
    7  = 9
    8  = 8
    9  = 7
    10 = 13
    11 = 14
    12 = 12
    13 = 11
    14 = 15
    15 = 10
    16 = 18
    17 = 21
    18 = 20
    19 = 17
    20 = 19
    21 = 16
    22 = 22
    23 = 23
    24 = 24
    25 = 25
}

# Snapshot original values for rows 2-25, columns A-F (1-6)
$original = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value()
    }
    $original[$r] = $rowVals
}

# Write back according to mapping
for ($destRow = 2; $destRow -le 25; $destRow++) {
    $srcRow = $mapping[$destRow]
    $vals = $original[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
